# Fix al borrar sport con equipos
# Rename a few API operation names in the "Diccionario" worksheet:
#   DELETE DeleteTeamByUserName  -> DELETE DeleteTeamByName
#   DELETE DeleteSportByUserName -> DELETE DeleteSportByName
#   GET GetSportByUserName       -> GET GetSportByName

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value  = "DELETE DeleteTeamByName"
$ws.Range("E10").Value = "DELETE DeleteSportByName"
$ws.Range("E28").Value = "GET GetSportByName"
